# Daily attendance processing - 2026-02-14 12:50:28 UTC
# Swap "Administrator, Miss Dina Nasr" -> "Miss Dina Nasr, Administrator"
# in the "Recorded By" column (G) wherever it currently appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "Administrator, Miss Dina Nasr"
$newValue = "Miss Dina Nasr, Administrator"

# Scan the used range of column G ("Recorded By") and flip every cell
# that still reads "Administrator, Miss Dina Nasr" to the reordered form.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
